$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.162124713418109
$ws.Range("C2").Value = 0.6102799498593754
$ws.Range("D2").Value = 0.5745248940500431
$ws.Range("E2").Value = 0.1959789392898088
$ws.Range("G2").Value = 0.00264742864926391
$ws.Range("I2").Value = 2.458167271194625
$ws.Range("J2").Value = 0.07560126133407863
$ws.Range("K2").Value = 1.769038056958266
$ws.Range("M2").Value = 0.6912718120050485
$ws.Range("N2").Value = 4.47231331304269
$ws.Range("B3").Value = 1.127532595743162
$ws.Range("C3").Value = 0.5979660578902894
$ws.Range("D3").Value = 0.5703747096809337
$ws.Range("E3").Value = 0.194994247533419
$ws.Range("G3").Value = 0.002652472904587884
$ws.Range("I3").Value = 2.442424863397079
$ws.Range("J3").Value = 0.07548774947014181
$ws.Range("K3").Value = 1.722387080737406
$ws.Range("M3").Value = 0.6794398837051716
$ws.Range("N3").Value = 4.460522937760828
$ws.Range("B4").Value = 1.107040958900626
$ws.Range("C4").Value = 0.5907390937300647
$ws.Range("D4").Value = 0.568096960883139
$ws.Range("E4").Value = 0.1944768266509378
$ws.Range("G4").Value = 0.002655732952265429
$ws.Range("I4").Value = 2.433524813918652
$ws.Range("J4").Value = 0.07544756848385248
$ws.Range("K4").Value = 1.69482752102266
$ws.Range("M4").Value = 0.6725640487035918
$ws.Range("N4").Value = 4.454047114754729
$ws.Range("B5").Value = 1.098878206224867
$ws.Range("C5").Value = 0.5878777411807334
$ws.Range("D5").Value = 0.5672367437112769
$ws.Range("E5").Value = 0.1942878899563212
$ws.Range("G5").Value = 0.002657102538213696
$ws.Range("I5").Value = 2.430090279845572
$ws.Range("J5").Value = 0.07543861634200155
$ws.Range("K5").Value = 1.683868834619346
$ws.Range("M5").Value = 0.6698597686181031
$ws.Range("N5").Value = 4.451599767902593
$ws.Range("B6").Value = 1.09753411735565
$ws.Range("C6").Value = 0.5874076635872996
$ws.Range("D6").Value = 0.5670980105391266
$ws.Range("E6").Value = 0.1942578406946147
$ws.Range("G6").Value = 0.002657332442407526
$ws.Range("I6").Value = 2.429531585256996
$ws.Range("J6").Value = 0.07543757811788154
$ws.Range("K6").Value = 1.68206556740472
$ws.Range("M6").Value = 0.669416620450825
$ws.Range("N6").Value = 4.451204948965682
$ws.Range("B7").Value = 1.106930113408595
$ws.Range("C7").Value = 0.5907001659784612
$ws.Range("D7").Value = 0.5680850844665883
$ws.Range("E7").Value = 0.1944741898531817
$ws.Range("G7").Value = 0.002655751256557509
$ws.Range("I7").Value = 2.433477716317014
$ws.Range("J7").Value = 0.07544741770087526
$ws.Range("K7").Value = 1.694678627589212
$ws.Range("M7").Value = 0.6725271825268493
$ws.Range("N7").Value = 4.454013333733499
$ws.Range("B8").Value = 1.150041862363878
$ws.Range("C8").Value = 0.6059646815511996
$ws.Range("D8").Value = 0.5730377501395481
$ws.Range("E8").Value = 0.1956213176093264
$ws.Range("G8").Value = 0.002649134191975339
$ws.Range("I8").Value = 2.452580162279517
$ws.Range("J8").Value = 0.07555599615478314
$ws.Range("K8").Value = 1.752727360097566
$ws.Range("M8").Value = 0.6871113297677098
$ws.Range("N8").Value = 4.468089281431531
$ws.Range("B9").Value = 1.240541762173876
$ws.Range("C9").Value = 0.6385606024073809
$ws.Range("D9").Value = 0.5848987635920651
$ws.Range("E9").Value = 0.1985632732808718
$ws.Range("G9").Value = 0.002637443983294372
$ws.Range("I9").Value = 2.496133012589397
$ws.Range("J9").Value = 0.07600321589249859
$ws.Range("K9").Value = 1.875200642420879
$ws.Range("M9").Value = 0.7188065256581098
$ws.Range("N9").Value = 4.501771262172582
$ws.Range("B10").Value = 1.31070448729696
$ws.Range("C10").Value = 0.6641541574429368
$ws.Range("D10").Value = 0.5949287013381479
$ws.Range("E10").Value = 0.2011483484733496
$ws.Range("G10").Value = 0.002629630176996688
$ws.Range("I10").Value = 2.531874232570885
$ws.Range("J10").Value = 0.07647491320045319
$ws.Range("K10").Value = 1.970511962620549
$ws.Range("M10").Value = 0.7439965303207217
$ws.Range("N10").Value = 4.530257327458656
$ws.Range("B11").Value = 1.343430230659578
$ws.Range("C11").Value = 0.6761597198528193
$ws.Range("D11").Value = 0.5997786515634402
$ws.Range("E11").Value = 0.2024167233769276
$ws.Range("G11").Value = 0.002626241858194065
$ws.Range("I11").Value = 2.548953222665006
$ws.Range("J11").Value = 0.07672065606275424
$ws.Range("K11").Value = 2.015043617248693
$ws.Range("M11").Value = 0.7558732878178915
$ws.Range("N11").Value = 4.544036200921113
$ws.Range("B12").Value = 1.355939460656543
$ws.Range("C12").Value = 0.680758458163524
$ws.Range("D12").Value = 0.6016565954921589
$ws.Range("E12").Value = 0.2029103345794923
$ws.Range("G12").Value = 0.00262498254774525
$ws.Range("I12").Value = 2.55553896672113
$ws.Range("J12").Value = 0.07681819837301873
$ws.Range("K12").Value = 2.032076387263999
$ws.Range("M12").Value = 0.7604310189881147
$ws.Range("N12").Value = 4.549372424515326
$ws.Range("B13").Value = 1.353240176383622
$ws.Range("C13").Value = 0.6797656980608622
$ws.Range("D13").Value = 0.6012503055608249
$ws.Range("E13").Value = 0.2028034345506384
$ws.Range("G13").Value = 0.002625252707692309
$ws.Range("I13").Value = 2.554115341469057
$ws.Range("J13").Value = 0.07679699137192131
$ws.Range("K13").Value = 2.028400520037962
$ws.Range("M13").Value = 0.7594467467203572
$ws.Range("N13").Value = 4.548217896731501
$ws.Range("B14").Value = 1.344457030995898
$ws.Range("C14").Value = 0.6765370070304755
$ws.Range("D14").Value = 0.5999323215947641
$ws.Range("E14").Value = 0.2024570663140679
$ws.Range("G14").Value = 0.002626137778224495
$ws.Range("I14").Value = 2.549492662722713
$ws.Range("J14").Value = 0.07672859102497753
$ws.Range("K14").Value = 2.016441509987771
$ws.Range("M14").Value = 0.7562470462396718
$ws.Range("N14").Value = 4.544472838245213
$ws.Range("B15").Value = 1.33909231263732
$ws.Range("C15").Value = 0.6745661857426342
$ws.Range("D15").Value = 0.5991304082585316
$ws.Range("E15").Value = 0.2022466390004567
$ws.Range("G15").Value = 0.002626683001846744
$ws.Range("I15").Value = 2.546676555884119
$ws.Range("J15").Value = 0.07668727799593
$ws.Range("K15").Value = 2.009138381025764
$ws.Range("M15").Value = 0.7542949909843273
$ws.Range("N15").Value = 4.542194323817654
$ws.Range("B16").Value = 1.308582088349112
$ws.Range("C16").Value = 0.6633768945457064
$ws.Range("D16").Value = 0.594617531768904
$ws.Range("E16").Value = 0.2010673181301286
$ws.Range("G16").Value = 0.002629854945245782
$ws.Range("I16").Value = 2.530774614708989
$ws.Range("J16").Value = 0.07645948072656594
$ws.Range("K16").Value = 1.967625404731223
$ws.Range("M16").Value = 0.7432287780865252
$ws.Range("N16").Value = 4.529373396126744
$ws.Range("B17").Value = 1.290072418563057
$ws.Range("C17").Value = 0.6566057951290531
$ws.Range("D17").Value = 0.5919226511128954
$ws.Range("E17").Value = 0.200367523140244
$ws.Range("G17").Value = 0.002631843310201315
$ws.Range("I17").Value = 2.52122960226076
$ws.Range("J17").Value = 0.07632771891403323
$ws.Range("K17").Value = 1.942459784699935
$ws.Range("M17").Value = 0.7365471351779362
$ws.Range("N17").Value = 4.521718645863018
$ws.Range("B18").Value = 1.27950220520259
$ws.Range("C18").Value = 0.6527453795578424
$ws.Range("D18").Value = 0.5903996633733186
$ws.Range("E18").Value = 0.1999737167028321
$ws.Range("G18").Value = 0.002633002618361191
$ws.Range("I18").Value = 2.515816735528958
$ws.Range("J18").Value = 0.07625486622374211
$ws.Range("K18").Value = 1.928095646162234
$ws.Range("M18").Value = 0.732743334500725
$ws.Range("N18").Value = 4.517393028079198
$ws.Range("B19").Value = 1.275936362043211
$ws.Range("C19").Value = 0.6514441638880442
$ws.Range("D19").Value = 0.5898886475659424
$ws.Range("E19").Value = 0.199841873863388
$ws.Range("G19").Value = 0.002633397832424271
$ws.Range("I19").Value = 2.513997275163817
$ws.Range("J19").Value = 0.07623070324305559
$ws.Range("K19").Value = 1.92325114419063
$ws.Range("M19").Value = 0.7314621774949188
$ws.Range("N19").Value = 4.515941691289953
$ws.Range("B20").Value = 1.292034929410875
$ws.Range("C20").Value = 0.6573230545041042
$ws.Range("D20").Value = 0.5922067271592368
$ws.Range("E20").Value = 0.2004411172752505
$ws.Range("G20").Value = 0.002631630026159754
$ws.Range("I20").Value = 2.522237695531331
$ws.Range("J20").Value = 0.076341441584308
$ws.Range("K20").Value = 1.945127270167177
$ws.Range("M20").Value = 0.7372543383001684
$ws.Range("N20").Value = 4.522525514565388
$ws.Range("B21").Value = 1.347033682739266
$ws.Range("C21").Value = 0.677483924766733
$ws.Range("D21").Value = 0.600318322066812
$ws.Range("E21").Value = 0.2025584418262412
$ws.Range("G21").Value = 0.002625877167428329
$ws.Range("I21").Value = 2.55084724188896
$ws.Range("J21").Value = 0.07674856012044273
$ws.Range("K21").Value = 2.019949552021217
$ws.Range("M21").Value = 0.7571852387740634
$ws.Range("N21").Value = 4.545569633712631
$ws.Range("B22").Value = 1.383658994300674
$ws.Range("C22").Value = 0.6909663305441427
$ws.Range("D22").Value = 0.6058609033509015
$ws.Range("E22").Value = 0.2040197968608126
$ws.Range("G22").Value = 0.002622255844827752
$ws.Range("I22").Value = 2.570234961425712
$ws.Range("J22").Value = 0.07704077834168288
$ws.Range("K22").Value = 2.069839162798644
$ws.Range("M22").Value = 0.7705625764412787
$ws.Range("N22").Value = 4.561321024709201
$ws.Range("B23").Value = 1.364048960065531
$ws.Range("C23").Value = 0.6837424054152734
$ws.Range("D23").Value = 0.6028806334611261
$ws.Range("E23").Value = 0.2032327415990736
$ws.Range("G23").Value = 0.002624175982879637
$ws.Range("I23").Value = 2.559824136627455
$ws.Range("J23").Value = 0.07688242274972268
$ws.Range("K23").Value = 2.043121391965258
$ws.Range("M23").Value = 0.7633906252270748
$ws.Range("N23").Value = 4.552850845885359
$ws.Range("B24").Value = 1.291147456380315
$ws.Range("C24").Value = 0.6569986807056694
$ws.Range("D24").Value = 0.5920782143560075
$ws.Range("E24").Value = 0.2004078188456226
$ws.Range("G24").Value = 0.002631726401346141
$ws.Range("I24").Value = 2.521781703934579
$ws.Range("J24").Value = 0.07633522853581809
$ws.Range("K24").Value = 1.94392097637575
$ws.Range("M24").Value = 0.7369344947751912
$ws.Range("N24").Value = 4.522160494995347
$ws.Range("B25").Value = 1.215417184364469
$ws.Range("C25").Value = 0.6294552419624324
$ws.Range("D25").Value = 0.5814594234021229
$ws.Range("E25").Value = 0.1976931153822008
$ws.Range("G25").Value = 0.002640469758555028
$ws.Range("I25").Value = 2.483695603157912
$ws.Range("J25").Value = 0.07585711898344627
$ws.Range("K25").Value = 1.875200642420879
$ws.Range("M25").Value = 0.7098990525810791
$ws.Range("N25").Value = 4.492005106474807
